$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.429.06'
$ws.Range('E2').Value = '  +0.53%  '

$ws.Range('D3').Value = '2.015.91'
$ws.Range('E3').Value = '  +0.59%  '

$ws.Range('E4').Value = '  -0.01%  '

$ws.Range('D5').Value = '''260.23'
$ws.Range('E5').Value = '  +5.63%  '

$ws.Range('D6').Value = '''0.619'
$ws.Range('E6').Value = '  -1.50%  '

$ws.Range('E7').Value = '  +0.00%  '

$ws.Range('D8').Value = '''56.26'
$ws.Range('E8').Value = '  -6.22%  '

$ws.Range('D9').Value = '''0.388'
$ws.Range('E9').Value = '  +1.13%  '

$ws.Range('E10').Value = '  -3.89%  '

$ws.Range('E11').Value = '  -2.03%  '

$ws.Range('D12').Value = '''14.34'

$ws.Range('D13').Value = '2.312.76'
$ws.Range('E13').Value = '  +0.66%  '

$ws.Range('D14').Value = '''0.808'
$ws.Range('E14').Value = '  -4.96%  '

$ws.Range('D15').Value = '''20.88'
$ws.Range('E15').Value = '  -7.79%  '

$ws.Range('E16').Value = '  -3.45%  '

$ws.Range('D17').Value = '2.030.08'
$ws.Range('E17').Value = '  +0.86%  '

$ws.Range('D18').Value = '37.272.35'
$ws.Range('E18').Value = '  +0.23%  '

$ws.Range('D19').Value = '''69.80'
$ws.Range('E19').Value = '  -0.86%  '

$ws.Range('E20').Value = '  -2.59%  '

$ws.Range('E21').Value = '  -0.04%  '

$ws.Range('D22').Value = '''228.52'

$ws.Range('D23').Value = '''2.68'
$ws.Range('E23').Value = '  +8.56%  '

$ws.Range('E24').Value = '  -0.03%  '

$ws.Range('D25').Value = '''2.33'
$ws.Range('E25').Value = '  -1.50%  '

$ws.Range('D26').Value = '''164.79'
$ws.Range('E26').Value = '  +0.50%  '

$ws.Range('E27').Value = '  -4.50%  '

$ws.Range('D28').Value = '''19.73'
$ws.Range('E28').Value = '  +0.14%  '

$ws.Range('E29').Value = '  -9.39%  '

$ws.Range('D30').Value = '''1.31'
$ws.Range('E30').Value = '  -2.70%  '

$ws.Range('E31').Value = '  -0.90%  '

$ws.Range('E32').Value = '  -3.19%  '

$ws.Range('D33').Value = '''0.0649'
$ws.Range('E33').Value = '  -1.34%  '

$ws.Range('D34').Value = '''4.56'
$ws.Range('E34').Value = '  +1.12%  '

$ws.Range('D35').Value = '''2.41'
$ws.Range('E35').Value = '  -0.18%  '

$ws.Range('D36').Value = '''1.83'
$ws.Range('E36').Value = '  +1.15%  '

$ws.Range('E37').Value = '  +0.04%  '

$ws.Range('D38').Value = '''3.33'
$ws.Range('E38').Value = '  +1.24%  '

$ws.Range('E39').Value = '  -4.03%  '

$ws.Range('E40').Value = '  +4.06%  '

$ws.Range('E41').Value = '  +3.72%  '

$ws.Range('D42').Value = '''0.0939'
$ws.Range('E42').Value = '  -4.15%  '

$ws.Range('D43').Value = '''0.0214'
$ws.Range('E43').Value = '  -0.81%  '

$ws.Range('D44').Value = '1.392.03'
$ws.Range('E44').Value = '  +0.79%  '

$ws.Range('D45').Value = '''90.32'
$ws.Range('E45').Value = '  -0.89%  '

$ws.Range('E46').Value = '  -5.81%  '

$ws.Range('E47').Value = '  -2.13%  '

$ws.Range('D48').Value = '''7.11'
$ws.Range('E48').Value = '  -4.90%  '

$ws.Range('D49').Value = '''2.90'
$ws.Range('E49').Value = '  +1.81%  '

$ws.Range('D50').Value = '2.203.27'
$ws.Range('E50').Value = '  +0.60%  '

$ws.Range('E51').Value = '  -4.88%  '
